# Apply scheduled-runner profit/price data refresh across all leve sheets
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2214.8
$ws.Range("J70").Value = 2714.8
$ws.Range("L70").Value = 8144.400000000001
$ws.Range("N70").Value = -8684.400000000001
$ws.Range("H73").Value = 2214.8
$ws.Range("J73").Value = 2714.8
$ws.Range("L73").Value = 8144.400000000001
$ws.Range("N73").Value = -10016.4
$ws.Range("H106").Value = 1500
$ws.Range("I106").Value = 1500
$ws.Range("K106").Value = 1500
$ws.Range("M106").Value = -869
$ws.Range("H132").Value = 62505544
$ws.Range("I132").Value = 71434264
$ws.Range("K132").Value = 214302792
$ws.Range("M132").Value = -214300262
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10754134
$ws.Range("I32").Value = 11765944
$ws.Range("K32").Value = 11765944
$ws.Range("M32").Value = -11765657
$ws.Range("H60").Value = 65013.5
$ws.Range("I60").Value = 65013.5
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 65013.5
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = -64280.5
$ws.Range("N60").ClearContents()
$ws.Range("H61").Value = 125002060
$ws.Range("I61").Value = 142858780
$ws.Range("K61").Value = 142858780
$ws.Range("M61").Value = -142858568
$ws.Range("H80").Value = 44993.332
$ws.Range("J80").Value = 44993.332
$ws.Range("L80").Value = 44993.332
$ws.Range("N80").Value = -46989.332
$ws.Range("H83").Value = 44993.332
$ws.Range("J83").Value = 44993.332
$ws.Range("L83").Value = 134979.996
$ws.Range("N83").Value = -144963.996
$ws.Range("H136").Value = 125002060
$ws.Range("I136").Value = 142858780
$ws.Range("K136").Value = 428576340
$ws.Range("M136").Value = -428573790

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 12757.8
$ws.Range("I26").Value = 10891.444
$ws.Range("K26").Value = 10891.444
$ws.Range("M26").Value = -10599.444
$ws.Range("H86").Value = 12730.909
$ws.Range("I86").Value = 14338
$ws.Range("K86").Value = 14338
$ws.Range("M86").Value = -13215
$ws.Range("H89").Value = 12730.909
$ws.Range("I89").Value = 14338
$ws.Range("K89").Value = 71690
$ws.Range("M89").Value = -66074
$ws.Range("H108").Value = 97748.75
$ws.Range("J108").Value = 94997.5
$ws.Range("L108").Value = 94997.5
$ws.Range("N108").Value = -102677.5
$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H131").Value = 30000
$ws.Range("J131").Value = 30000
$ws.Range("L131").Value = 30000
$ws.Range("N131").Value = -40080

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 359.14285
$ws.Range("I11").Value = 263
$ws.Range("K11").Value = 263
$ws.Range("M11").Value = -123
$ws.Range("H31").Value = 41671130
$ws.Range("I31").Value = 2537.318
$ws.Range("J31").Value = 156259740
$ws.Range("K31").Value = 2537.318
$ws.Range("L31").Value = 156259740
$ws.Range("M31").Value = -2242.318
$ws.Range("N31").Value = -156260330
$ws.Range("H34").Value = 41671130
$ws.Range("I34").Value = 2537.318
$ws.Range("J34").Value = 156259740
$ws.Range("K34").Value = 2537.318
$ws.Range("L34").Value = 156259740
$ws.Range("M34").Value = -2335.318
$ws.Range("N34").Value = -156260144
$ws.Range("H56").Value = 13000
$ws.Range("I56").Value = 6000
$ws.Range("K56").Value = 6000
$ws.Range("M56").Value = -5155
$ws.Range("H108").Value = 26666.666
$ws.Range("J108").Value = 26666.666
$ws.Range("L108").Value = 26666.666
$ws.Range("N108").Value = -34346.666
$ws.Range("H122").Value = 2394206.8
$ws.Range("J122").Value = 5265393
$ws.Range("L122").Value = 15796179
$ws.Range("N122").Value = -15801079

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 26
$ws.Range("I6").Value = 26
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 78
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 35
$ws.Range("N6").ClearContents()
$ws.Range("H11").Value = 115.25
$ws.Range("I11").Value = 87
$ws.Range("K11").Value = 261
$ws.Range("M11").Value = -121
$ws.Range("H14").Value = 50025
$ws.Range("I14").Value = 50025
$ws.Range("K14").Value = 150075
$ws.Range("M14").Value = -149902
$ws.Range("H109").Value = 4760.1816
$ws.Range("I109").Value = 4836.3
$ws.Range("J109").Value = 3999
$ws.Range("K109").Value = 14508.9
$ws.Range("L109").Value = 11997
$ws.Range("M109").Value = -13468.9
$ws.Range("N109").Value = -14077
$ws.Range("H121").Value = 7143580.5
$ws.Range("I121").Value = 20000438
$ws.Range("K121").Value = 60001314
$ws.Range("M121").Value = -60000004
$ws.Range("H129").Value = 46668030
$ws.Range("J129").Value = 22223974
$ws.Range("L129").Value = 66671922
$ws.Range("N129").Value = -66681922

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4290.5
$ws.Range("I102").Value = 3250.8
$ws.Range("J102").Value = 5156.9165
$ws.Range("K102").Value = 3250.8
$ws.Range("L102").Value = 5156.9165
$ws.Range("M102").Value = -1628.8
$ws.Range("N102").Value = -8400.916499999999
$ws.Range("I122").Value = 3336
$ws.Range("J122").Value = 125003470
$ws.Range("K122").Value = 10008
$ws.Range("L122").Value = 375010410
$ws.Range("M122").Value = -7558
$ws.Range("N122").Value = -375015310

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 478.75
$ws.Range("I55").Value = 439.63635
$ws.Range("K55").Value = 439.63635
$ws.Range("M55").Value = -266.63635
$ws.Range("H136").Value = 1824797.1
$ws.Range("I136").Value = 1824797.1
$ws.Range("K136").Value = 5474391.300000001
$ws.Range("M136").Value = -5471841.300000001
$ws.Range("H139").Value = 65781.2
$ws.Range("I139").Value = 60324.5
$ws.Range("J139").Value = 69419
$ws.Range("K139").Value = 60324.5
$ws.Range("L139").Value = 69419
$ws.Range("M139").Value = -55184.5
$ws.Range("N139").Value = -79699

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H141").Value = 110998.75
$ws.Range("J141").Value = 111665
$ws.Range("L141").Value = 111665
$ws.Range("N141").Value = -122025

